# Ran code for averaged intensities on spiral schemes.
#
# The underlying notebook was re-run after adding three new "spiral" sampling
# schemes and re-positioning "Gaussian-Quadrature" right after the "Ring
# Perpendicular to *" group. That shifts the existing scheme rows down and
# moves the three "HexGrid-*" rows to the bottom of the table, with the
# averaged-intensity values (always 1) recomputed for every row/column pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final ordering of scheme names for rows 10-19 (rows 3-9 are unchanged).
$schemeNames = @(
    "Gaussian-Quadrature",
    "Spiral-90deg-10rot-5space",
    "Spiral-90deg-15rot-5space",
    "Spiral-90deg-10rot-3space",
    "NoRotation-tilt60deg",
    "Rotation-NoTilt",
    "Rotation-60detTilt",
    "HexGrid-90degTilt5degRes",
    "HexGrid-90degTilt22p5degRes",
    "HexGrid-60degTilt5degRes"
)

$startRow = 10
$lastExistingRow = 16

for ($i = 0; $i -lt $schemeNames.Length; $i++) {
    $row = $startRow + $i

    if ($row -gt $lastExistingRow) {
        # Brand new row: clone the formatting (bold index column / plain data
        # columns) from the last existing data row before filling it in.
        $ws.Cells.Item($lastExistingRow, 1).Copy()
        $ws.Cells.Item($row, 1).PasteSpecial(-4122)
        $excel.CutCopyMode = $false
    }

    # Column A: sequential row index (0-based, same scheme as before)
    $ws.Cells.Item($row, 1).Value = $row - 2

    # Column B: the (possibly reordered / new) scheme name
    $ws.Cells.Item($row, 2).Value = $schemeNames[$i]

    # Columns C-M: averaged intensity values, all 1 for every HKL / pairing column
    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}
